$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.302.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.933.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7490"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.88"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3180"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07145"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7819"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08044"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.930.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.400"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.54"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.308.94"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.080"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007975"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.179.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.677"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.555"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.08"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1300"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.197"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.371"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.544"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.420"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.156"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.335"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.58%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05259"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7587"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.787"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01953"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.797"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.498"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4524"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.981"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8403"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.992"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.678"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.80"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.61"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1219"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.41%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "958.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.93%  "
